$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36:57 down to 37:58
$ws.Rows("36:36").Insert()

# Populate the new row 36 with the linkeR entry
$ws.Cells.Item(36, 1).Value = "linkeR: Effortless Linked Views for Shiny Applications"
$ws.Cells.Item(36, 2).Value = "linkeR makes it effortless to create linked views in Shiny applications. When users interact with one component (like clicking a map marker), all related components (tables, charts, other maps) automatically update to show corresponding information."
$ws.Cells.Item(36, 3).Value = "Jake Wagoner"
$ws.Cells.Item(36, 4).Value = "jakew@sci.utah.edu"
$ws.Cells.Item(36, 5).Value = "Yes"
$ws.Cells.Item(36, 7).Value = "Development"
$ws.Cells.Item(36, 8).Value = "MIT"
$ws.Cells.Item(36, 9).Value = "R"
$ws.Cells.Item(36, 10).Value = "R Shiny Developers"
$ws.Cells.Item(36, 11).Value = "Moderate Programming"
$ws.Cells.Item(36, 12).Value = "Developer Tool"
$ws.Cells.Item(36, 14).Value = "https://epiforesite.github.io/linkeR/"
$ws.Cells.Item(36, 15).Value = "https://github.com/EpiForeSITE/linkeR/"
